$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    2  = 1
    3  = 2
    4  = -1
    5  = 5
    15 = 3
    26 = 5
    29 = 0
    31 = 0
    40 = 0
    44 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
